$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------
# Helper-ish constants (Excel COM enum values)
# xlPasteFormats = -4122, xlPasteValues = -4163
# -------------------------------------------------------------------

# ---------------------------------------------------------------
# 1) Ativação: date text "01/01/2012" -> "01/01/2023"
#    (this shared string is reused verbatim by B8/C8 and B15/C15).
#    We route the new text through a scratch formula cell first so the
#    engine stores it as literal text instead of re-parsing it as a
#    date serial number.
# ---------------------------------------------------------------
$ws.Range("Z1").Formula = "=""01/01/2023"""
$ws.Range("Z1").Copy()
$ws.Range("B8").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("C8").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# ---------------------------------------------------------------
# 2) Row 11 "Objectives:" gains B11/C11 with the English objectives text.
#    Pull formatting from row 13 (same B/C column styles) first, since a
#    brand-new cell otherwise inherits the row's style instead of the
#    column's.
# ---------------------------------------------------------------
$objectivesText = "Rheology is the science that studies the flow of materials. Your knowledge is necessary to understand the processes of forming materials. The course aims to provide students with the basic and applied concepts of rheology and familiarize them with experimental methods for evaluating the rheological properties of materials."

$ws.Range("B13").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("B11").Value = $objectivesText
$ws.Range("C11").Value = $objectivesText

# ---------------------------------------------------------------
# 3) Row 14 "Short syllabus:" gains B14/C14 with the English short syllabus.
# ---------------------------------------------------------------
$shortSyllabusText = "Flow of Newtonian and non-Newtonian fluids. Viscosity and rheometry. viscoelasticity. Applications."

$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("B14").Value = $shortSyllabusText
$ws.Range("C14").Value = $shortSyllabusText

# ---------------------------------------------------------------
# 4) Row 16 "Syllabus:" gains B16/C16 with the English full syllabus text.
# ---------------------------------------------------------------
$syllabusText = "1. Introduction. 2. Stress and deformation. 3. Types of deformation and flow of materials. 4. Fundamental equations of rheology. Flow of Newtonian and non-Newtonian fluids. 5. Viscosimetry and rheometry. 6. Rheology of dispersed systems. Colloids and emulsions. diluted solutions. Capillary viscosimetry. 7. Rheology of molten polymers. 8. Viscoelasticity. 9. Dynamic-mechanical behavior of materials. 10. Applications."

$ws.Range("B21").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C21").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("B16").Value = $syllabusText
$ws.Range("C16").Value = $syllabusText
